# Industry to ISIC Code Map.xlsx
# Commit: "Use US EPS values for bldgs/ICpUEfEBE, indst/CoNEPPpCAPS, indst/EoDfIP,
#          indst/ItICM, and trans/RTMF"
#
# For this workbook, the substantive change is on the "About" sheet: a new
# note is inserted ("The EU EPS uses values from the US EPS.") right after the
# "Notes" heading block, pushing the rest of the notes down by two rows (one
# row for the new text, one blank spacer row, matching the blank-row style
# used throughout the rest of the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Make room for the new note: insert two rows above the old row 11
# ("Each industry doesn't need to be mapped ..."), shifting everything from
# the old row 11 onward down to row 13 onward.
$ws.Rows("11:12").Insert()

# Write the new note into the freshly inserted row 11 (row 12 is left blank,
# consistent with the other section-separating blank rows on this sheet).
$ws.Range("A11").Value = "The EU EPS uses values from the US EPS."

# Leave the sheet's cell selection where the author left it when saving.
$ws.Range("A12").Select() | Out-Null
